$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated sval data (filtered save games) for rows 2-7, columns B-G
$data = @{
    2 = @(0.127881588408715, 0.3127903958511391, 0.1575252929769615, 0.496779210170732, 0, 1.094976487407548)
    3 = @(3.230985683306322, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 1, 6.201049113329182)
    4 = @(0.04763786555579896, 1.667794583268128, 337.1190423067083, 8.660232485948974, 0, 347.4947072414812)
    5 = @(0.3048080303191223, 0.3127903958511391, 0.8054896365839992, 8.660232485948974, 1, 10.08332054870323)
    6 = @(0.6753301551942219, 0.00007097389502863649, 0.1575252929769615, 0.496779210170732, 0, 1.329705632236944)
    7 = @(3.230985683306322, 10.29869402782916, 0.8054896365839992, 645.3272768299601, 1, 659.6624461776795)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $col = 2 + $i  # B=2 ... G=7
        $ws.Cells.Item($row, $col).Value = $vals[$i]
    }
}

$wb.Save()
